$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Wrong_Entity_NonEvent_as_Event"
$ws.Range("A3").Value = "Correct"
$ws.Range("A4").Value = "Wrong_Tag_E_as_I"
$ws.Range("A5").Value = "Wrong_Entity_Event_as_NonEvent"
